$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B17: change from inline string "3" to numeric 3
$ws.Range("B17").Value = 3

# Add new row 18
$ws.Range("A18").Value = "Ying Tang"
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "4"
$ws.Range("C18").Value = " show through several experiments "
$ws.Range("D18").Value = "FBK"
$ws.Range("E18").Value = "RES"
$ws.Range("F18").Value = "9cb2103f-10a8-4188-b35f-b6e342d90889"
$ws.Range("G18").Value = "rJwelMbR-_annotated.xlsx"
$ws.Range("H18").Value = "The authors show through several experiments that the divide and conquer (DnC) technique can solve more complex tasks than can be solved with conventional policy gradient methods (TRPO is used as the baseline)."
